# Updated cryptos list - refresh Price / Volume(1h) values, and restore the
# natural coinranking.com ordering for rows whose Coin/Link/Price/Volume cells
# got shuffled (XRP/BNB, Toncoin/WEMIXToken, Kaspa/EthereumClassic,
# TerraClassic/FraxShare, Cronos/Aave/TrustWalletToken) plus one row that now
# shows a different coin entirely (HuobiToken -> Celestia).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($Row, $Coin, $Link, $Price, $Volume) {
    if ($Coin) { $ws.Cells.Item($Row, 2).Value = $Coin }
    if ($Link) { $ws.Cells.Item($Row, 3).Value = $Link }
    $ws.Cells.Item($Row, 4).NumberFormat = "@"
    $ws.Cells.Item($Row, 4).Value = $Price
    $ws.Cells.Item($Row, 4).Style = "Normal"
    $ws.Cells.Item($Row, 5).NumberFormat = "@"
    $ws.Cells.Item($Row, 5).Value = $Volume
    $ws.Cells.Item($Row, 5).Style = "Normal"
}

Set-Row 2 $null $null "43.849.70" "  -0.63%  "
Set-Row 3 $null $null "2.241.90" "  -2.96%  "
Set-Row 4 $null $null "1.00" "  +0.44%  "
Set-Row 5 "BNB" "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb" "230.85" "  -0.70%  "
Set-Row 6 "XRP" "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp" "0.645" "  +3.68%  "
Set-Row 7 $null $null "63.37" "  -0.04%  "
Set-Row 8 $null $null "1.00" "  +0.11%  "
Set-Row 9 $null $null "0.456" "  +6.06%  "
Set-Row 10 $null $null "0.0970" "  +0.80%  "
Set-Row 11 $null $null "57.16" "  -1.74%  "
Set-Row 12 $null $null "26.20" "  -1.26%  "
Set-Row 13 $null $null "0.106" "  +1.11%  "
Set-Row 14 $null $null "2.577.92" "  -2.25%  "
Set-Row 15 $null $null "15.46" "  -4.48%  "
Set-Row 16 $null $null "6.12" "  +2.16%  "
Set-Row 17 $null $null "0.836" "  +1.18%  "
Set-Row 18 $null $null "2.245.35" "  -2.01%  "
Set-Row 19 $null $null "43.773.28" "  -0.73%  "
Set-Row 20 $null $null "0.0₃0987" "  +3.10%  "
Set-Row 21 $null $null "72.78" "  -1.99%  "
Set-Row 22 $null $null "6.07" "  -3.30%  "
Set-Row 23 $null $null "248.95" "  -2.43%  "
Set-Row 24 $null $null "1.00" "  +0.06%  "
Set-Row 25 $null $null "2.41" "  -6.69%  "
Set-Row 26 "Toncoin" "https://coinranking.com/coin/67YlI0K1b+toncoin-ton" "2.29" "  -3.05%  "
Set-Row 27 "WEMIXToken" "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix" "3.31" "  +16.41%  "
Set-Row 28 $null $null "9.87" "  -1.01%  "
Set-Row 29 $null $null "171.06" "  -0.11%  "
Set-Row 30 "EthereumClassic" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc" "20.96" "  +0.77%  "
Set-Row 31 "Kaspa" "https://coinranking.com/coin/V8GxkwWow+kaspa-kas" "0.138" "  -2.73%  "
Set-Row 32 $null $null "1.41" "  -4.10%  "
Set-Row 33 $null $null "0.125" "  +2.20%  "
Set-Row 34 $null $null "0.0699" "  -0.97%  "
Set-Row 35 $null $null "4.79" "  -0.33%  "
Set-Row 36 $null $null "4.93" "  -4.15%  "
Set-Row 37 $null $null "3.65" "  -4.05%  "
Set-Row 38 $null $null "6.43" "  -3.38%  "
Set-Row 39 $null $null "2.28" "  -6.37%  "
Set-Row 40 $null $null "0.0258" "  +2.41%  "
Set-Row 41 $null $null "1.00" "  +0.16%  "
Set-Row 42 "FraxShare" "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs" "8.29" "  -5.88%  "
Set-Row 43 "TerraClassic" "https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc" "0.000217" "  -4.77%  "
Set-Row 44 $null $null "17.07" "  -1.32%  "
Set-Row 45 "Aave" "https://coinranking.com/coin/ixgUfzmLR+aave-aave" "97.32" "  -2.90%  "
Set-Row 46 "TrustWalletToken" "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt" "1.19" "  -3.18%  "
Set-Row 47 "Cronos" "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro" "0.0945" "  -2.82%  "
Set-Row 48 $null $null "4.36" "  -5.76%  "
Set-Row 49 $null $null "2.32" "  -0.25%  "
Set-Row 50 $null $null "1.427.54" "  -4.17%  "
Set-Row 51 "Celestia" "https://coinranking.com/coin/YQcD0lBl7+celestia-tia" "9.83" "  -8.17%  "
